# Add a new "TR200/M" optical-post line item to the Main parts list,
# just above the existing TBS1000C (Tektronix oscilloscope) row — i.e.
# insert a brand-new row 158 and push everything from the old row 158
# downward by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 158 (same formatting/formulas as the row above
# are inherited automatically, and every formula below — the shared
# F-column multiplications as well as the grand-total SUM — gets its
# references shifted down by Excel automatically).
$ws.Rows("158:158").Insert()

# Populate the new row with the TR200/M post used for the excitation
# cage alignment jig.
$ws.Range("A158").Value = "TR200/M"
$ws.Range("B158").Value = "Thorlabs"
$ws.Range("C158").Value = "Ø12.7 mm Optical Post, SS, M4 Setscrew, M6 Tap, L = 200 mm"
$ws.Range("D158").Value = 2
$ws.Range("E158").Value = 9
$ws.Range("F158").Formula = "=E158*D158"
$ws.Range("I158").Value = "Excitation cage alignment jig, see details at https://github.com/mesoSPIM/benchtop-hardware/wiki/Excitation_alignment"

# Match the author's final on-screen selection/scroll position.
$ws.Range("A158").Select()
$excel.ActiveWindow.ScrollRow = 117
